# Update the worksheet date and the 100 arithmetic problem/answer cells
# to the values from the new revision. Every "old" text below is a
# unique, exact run of text in the document (the date line, or a single
# table-cell run like "24+14=38"), so a simple Find/Replace across the
# whole document body is sufficient and safe (no ambiguous partial
# matches exist between any two target strings).
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-01 Saturday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-02-02 Sunday", 2) | Out-Null
$d.Content.Find.Execute("24+14=38", $true, $true, $false, $false, $false, $true, 1, $false, "68-43=25", 2) | Out-Null
$d.Content.Find.Execute("35+47=82", $true, $true, $false, $false, $false, $true, 1, $false, "80-38=42", 2) | Out-Null
$d.Content.Find.Execute("28+7=35", $true, $true, $false, $false, $false, $true, 1, $false, "0+94=94", 2) | Out-Null
$d.Content.Find.Execute("56+23=79", $true, $true, $false, $false, $false, $true, 1, $false, "49+12=61", 2) | Out-Null
$d.Content.Find.Execute("86-43=43", $true, $true, $false, $false, $false, $true, 1, $false, "36+46=82", 2) | Out-Null
$d.Content.Find.Execute("38+50=88", $true, $true, $false, $false, $false, $true, 1, $false, "91-5=86", 2) | Out-Null
$d.Content.Find.Execute("50+0=50", $true, $true, $false, $false, $false, $true, 1, $false, "26+6=32", 2) | Out-Null
$d.Content.Find.Execute("17+20=37", $true, $true, $false, $false, $false, $true, 1, $false, "2+7=9", 2) | Out-Null
$d.Content.Find.Execute("80-2=78", $true, $true, $false, $false, $false, $true, 1, $false, "40-14=26", 2) | Out-Null
$d.Content.Find.Execute("76-30=46", $true, $true, $false, $false, $false, $true, 1, $false, "9+60=69", 2) | Out-Null
$d.Content.Find.Execute("1+11=12", $true, $true, $false, $false, $false, $true, 1, $false, "18+8=26", 2) | Out-Null
$d.Content.Find.Execute("48-4=44", $true, $true, $false, $false, $false, $true, 1, $false, "68-32=36", 2) | Out-Null
$d.Content.Find.Execute("9+86=95", $true, $true, $false, $false, $false, $true, 1, $false, "15+55=70", 2) | Out-Null
$d.Content.Find.Execute("97-23=74", $true, $true, $false, $false, $false, $true, 1, $false, "12+38=50", 2) | Out-Null
$d.Content.Find.Execute("30+55=85", $true, $true, $false, $false, $false, $true, 1, $false, "83-7=76", 2) | Out-Null
$d.Content.Find.Execute("53-38=15", $true, $true, $false, $false, $false, $true, 1, $false, "72-43=29", 2) | Out-Null
$d.Content.Find.Execute("93-46=47", $true, $true, $false, $false, $false, $true, 1, $false, "96-58=38", 2) | Out-Null
$d.Content.Find.Execute("71-50=21", $true, $true, $false, $false, $false, $true, 1, $false, "13+79=92", 2) | Out-Null
$d.Content.Find.Execute("2+64=66", $true, $true, $false, $false, $false, $true, 1, $false, "87-26=61", 2) | Out-Null
$d.Content.Find.Execute("77-62=15", $true, $true, $false, $false, $false, $true, 1, $false, "20+35=55", 2) | Out-Null
$d.Content.Find.Execute("34+19=53", $true, $true, $false, $false, $false, $true, 1, $false, "53-19=34", 2) | Out-Null
$d.Content.Find.Execute("71-28=43", $true, $true, $false, $false, $false, $true, 1, $false, "2+9=11", 2) | Out-Null
$d.Content.Find.Execute("8+32=40", $true, $true, $false, $false, $false, $true, 1, $false, "8+64=72", 2) | Out-Null
$d.Content.Find.Execute("38+28=66", $true, $true, $false, $false, $false, $true, 1, $false, "5+56=61", 2) | Out-Null
$d.Content.Find.Execute("16+36=52", $true, $true, $false, $false, $false, $true, 1, $false, "1+50=51", 2) | Out-Null
$d.Content.Find.Execute("98-43=55", $true, $true, $false, $false, $false, $true, 1, $false, "7+54=61", 2) | Out-Null
$d.Content.Find.Execute("16-15=1", $true, $true, $false, $false, $false, $true, 1, $false, "90-14=76", 2) | Out-Null
$d.Content.Find.Execute("92-27=65", $true, $true, $false, $false, $false, $true, 1, $false, "17+17=34", 2) | Out-Null
$d.Content.Find.Execute("79+18=97", $true, $true, $false, $false, $false, $true, 1, $false, "19-5=14", 2) | Out-Null
$d.Content.Find.Execute("37-27=10", $true, $true, $false, $false, $false, $true, 1, $false, "93-79=14", 2) | Out-Null
$d.Content.Find.Execute("59-37=22", $true, $true, $false, $false, $false, $true, 1, $false, "95-23=72", 2) | Out-Null
$d.Content.Find.Execute("49-27=22", $true, $true, $false, $false, $false, $true, 1, $false, "17+57=74", 2) | Out-Null
$d.Content.Find.Execute("94-14=80", $true, $true, $false, $false, $false, $true, 1, $false, "47-34=13", 2) | Out-Null
$d.Content.Find.Execute("23+20=43", $true, $true, $false, $false, $false, $true, 1, $false, "2+4=6", 2) | Out-Null
$d.Content.Find.Execute("87+8=95", $true, $true, $false, $false, $false, $true, 1, $false, "66+20=86", 2) | Out-Null
$d.Content.Find.Execute("52+45=97", $true, $true, $false, $false, $false, $true, 1, $false, "38+56=94", 2) | Out-Null
$d.Content.Find.Execute("64-48=16", $true, $true, $false, $false, $false, $true, 1, $false, "47+47=94", 2) | Out-Null
$d.Content.Find.Execute("67-30=37", $true, $true, $false, $false, $false, $true, 1, $false, "45+41=86", 2) | Out-Null
$d.Content.Find.Execute("87-66=21", $true, $true, $false, $false, $false, $true, 1, $false, "22+43=65", 2) | Out-Null
$d.Content.Find.Execute("34+65=99", $true, $true, $false, $false, $false, $true, 1, $false, "38+11=49", 2) | Out-Null
$d.Content.Find.Execute("29+1=30", $true, $true, $false, $false, $false, $true, 1, $false, "56-18=38", 2) | Out-Null
$d.Content.Find.Execute("83+4=87", $true, $true, $false, $false, $false, $true, 1, $false, "56+41=97", 2) | Out-Null
$d.Content.Find.Execute("42-18=24", $true, $true, $false, $false, $false, $true, 1, $false, "17+42=59", 2) | Out-Null
$d.Content.Find.Execute("77-9=68", $true, $true, $false, $false, $false, $true, 1, $false, "20+46=66", 2) | Out-Null
$d.Content.Find.Execute("36+20=56", $true, $true, $false, $false, $false, $true, 1, $false, "52-31=21", 2) | Out-Null
$d.Content.Find.Execute("31+62=93", $true, $true, $false, $false, $false, $true, 1, $false, "28+26=54", 2) | Out-Null
$d.Content.Find.Execute("4+64=68", $true, $true, $false, $false, $false, $true, 1, $false, "71-53=18", 2) | Out-Null
$d.Content.Find.Execute("69-21=48", $true, $true, $false, $false, $false, $true, 1, $false, "52-41=11", 2) | Out-Null
$d.Content.Find.Execute("80+13=93", $true, $true, $false, $false, $false, $true, 1, $false, "5+81=86", 2) | Out-Null
$d.Content.Find.Execute("68+18=86", $true, $true, $false, $false, $false, $true, 1, $false, "53-0=53", 2) | Out-Null
$d.Content.Find.Execute("60-34=26", $true, $true, $false, $false, $false, $true, 1, $false, "14+25=39", 2) | Out-Null
$d.Content.Find.Execute("92-45=47", $true, $true, $false, $false, $false, $true, 1, $false, "88-17=71", 2) | Out-Null
$d.Content.Find.Execute("95-76=19", $true, $true, $false, $false, $false, $true, 1, $false, "39+19=58", 2) | Out-Null
$d.Content.Find.Execute("59-41=18", $true, $true, $false, $false, $false, $true, 1, $false, "84-72=12", 2) | Out-Null
$d.Content.Find.Execute("4+24=28", $true, $true, $false, $false, $false, $true, 1, $false, "15+81=96", 2) | Out-Null
$d.Content.Find.Execute("72-36=36", $true, $true, $false, $false, $false, $true, 1, $false, "14+59=73", 2) | Out-Null
$d.Content.Find.Execute("21-13=8", $true, $true, $false, $false, $false, $true, 1, $false, "94-25=69", 2) | Out-Null
$d.Content.Find.Execute("56-34=22", $true, $true, $false, $false, $false, $true, 1, $false, "52-52=0", 2) | Out-Null
$d.Content.Find.Execute("84+7=91", $true, $true, $false, $false, $false, $true, 1, $false, "37+2=39", 2) | Out-Null
$d.Content.Find.Execute("91-40=51", $true, $true, $false, $false, $false, $true, 1, $false, "79-77=2", 2) | Out-Null
$d.Content.Find.Execute("80-28=52", $true, $true, $false, $false, $false, $true, 1, $false, "45+5=50", 2) | Out-Null
$d.Content.Find.Execute("85-40=45", $true, $true, $false, $false, $false, $true, 1, $false, "14+38=52", 2) | Out-Null
$d.Content.Find.Execute("28+8=36", $true, $true, $false, $false, $false, $true, 1, $false, "18-2=16", 2) | Out-Null
$d.Content.Find.Execute("88-36=52", $true, $true, $false, $false, $false, $true, 1, $false, "73-33=40", 2) | Out-Null
$d.Content.Find.Execute("3+37=40", $true, $true, $false, $false, $false, $true, 1, $false, "38-20=18", 2) | Out-Null
$d.Content.Find.Execute("79-26=53", $true, $true, $false, $false, $false, $true, 1, $false, "24+40=64", 2) | Out-Null
$d.Content.Find.Execute("80-65=15", $true, $true, $false, $false, $false, $true, 1, $false, "67-2=65", 2) | Out-Null
$d.Content.Find.Execute("9+49=58", $true, $true, $false, $false, $false, $true, 1, $false, "68+8=76", 2) | Out-Null
$d.Content.Find.Execute("28-1=27", $true, $true, $false, $false, $false, $true, 1, $false, "50+4=54", 2) | Out-Null
$d.Content.Find.Execute("75-60=15", $true, $true, $false, $false, $false, $true, 1, $false, "27+26=53", 2) | Out-Null
$d.Content.Find.Execute("44+28=72", $true, $true, $false, $false, $false, $true, 1, $false, "63-7=56", 2) | Out-Null
$d.Content.Find.Execute("33+31=64", $true, $true, $false, $false, $false, $true, 1, $false, "48+13=61", 2) | Out-Null
$d.Content.Find.Execute("29+20=49", $true, $true, $false, $false, $false, $true, 1, $false, "10+33=43", 2) | Out-Null
$d.Content.Find.Execute("14+61=75", $true, $true, $false, $false, $false, $true, 1, $false, "78-36=42", 2) | Out-Null
$d.Content.Find.Execute("46-13=33", $true, $true, $false, $false, $false, $true, 1, $false, "89-23=66", 2) | Out-Null
$d.Content.Find.Execute("45+32=77", $true, $true, $false, $false, $false, $true, 1, $false, "5+71=76", 2) | Out-Null
$d.Content.Find.Execute("61+0=61", $true, $true, $false, $false, $false, $true, 1, $false, "39+31=70", 2) | Out-Null
$d.Content.Find.Execute("85-38=47", $true, $true, $false, $false, $false, $true, 1, $false, "91-11=80", 2) | Out-Null
$d.Content.Find.Execute("28+6=34", $true, $true, $false, $false, $false, $true, 1, $false, "42-37=5", 2) | Out-Null
$d.Content.Find.Execute("49+1=50", $true, $true, $false, $false, $false, $true, 1, $false, "23+34=57", 2) | Out-Null
$d.Content.Find.Execute("79-47=32", $true, $true, $false, $false, $false, $true, 1, $false, "14+8=22", 2) | Out-Null
$d.Content.Find.Execute("98-50=48", $true, $true, $false, $false, $false, $true, 1, $false, "94-37=57", 2) | Out-Null
$d.Content.Find.Execute("68+5=73", $true, $true, $false, $false, $false, $true, 1, $false, "99-45=54", 2) | Out-Null
$d.Content.Find.Execute("51+37=88", $true, $true, $false, $false, $false, $true, 1, $false, "94-61=33", 2) | Out-Null
$d.Content.Find.Execute("74-63=11", $true, $true, $false, $false, $false, $true, 1, $false, "72-57=15", 2) | Out-Null
$d.Content.Find.Execute("76-6=70", $true, $true, $false, $false, $false, $true, 1, $false, "51-47=4", 2) | Out-Null
$d.Content.Find.Execute("55+7=62", $true, $true, $false, $false, $false, $true, 1, $false, "29+39=68", 2) | Out-Null
$d.Content.Find.Execute("76-67=9", $true, $true, $false, $false, $false, $true, 1, $false, "96-36=60", 2) | Out-Null
$d.Content.Find.Execute("21-9=12", $true, $true, $false, $false, $false, $true, 1, $false, "3+4=7", 2) | Out-Null
$d.Content.Find.Execute("28+13=41", $true, $true, $false, $false, $false, $true, 1, $false, "41+21=62", 2) | Out-Null
$d.Content.Find.Execute("20+29=49", $true, $true, $false, $false, $false, $true, 1, $false, "81-49=32", 2) | Out-Null
$d.Content.Find.Execute("19+14=33", $true, $true, $false, $false, $false, $true, 1, $false, "12+71=83", 2) | Out-Null
$d.Content.Find.Execute("30-15=15", $true, $true, $false, $false, $false, $true, 1, $false, "74+5=79", 2) | Out-Null
$d.Content.Find.Execute("27+21=48", $true, $true, $false, $false, $false, $true, 1, $false, "50-22=28", 2) | Out-Null
$d.Content.Find.Execute("67-29=38", $true, $true, $false, $false, $false, $true, 1, $false, "92-37=55", 2) | Out-Null
$d.Content.Find.Execute("85-66=19", $true, $true, $false, $false, $false, $true, 1, $false, "74-32=42", 2) | Out-Null
$d.Content.Find.Execute("98-42=56", $true, $true, $false, $false, $false, $true, 1, $false, "63+17=80", 2) | Out-Null
$d.Content.Find.Execute("38+51=89", $true, $true, $false, $false, $false, $true, 1, $false, "12+52=64", 2) | Out-Null
$d.Content.Find.Execute("18+80=98", $true, $true, $false, $false, $false, $true, 1, $false, "16-0=16", 2) | Out-Null
$d.Content.Find.Execute("89-88=1", $true, $true, $false, $false, $false, $true, 1, $false, "93-18=75", 2) | Out-Null
